# Apply crypto price/volume update, generated from commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '70.886.07'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +6.07%  '

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '3.651.70'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +17.67%  '

# Row 4
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '617.82'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +7.14%  '

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '181.36'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +2.19%  '

# Row 7
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '3.648.62'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  +17.70%  '

# Row 8
$ws.Range("E8").Value = '  +0.03%  '

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.542'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +5.49%  '

# Row 10
$ws.Range("E10").Value = '  +8.25%  '

# Row 11
$ws.Range("E11").Value = '  +5.15%  '

# Row 12
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.503'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +7.58%  '

# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '40.42'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +11.87%  '

# Row 14
$ws.Range("E14").Value = '  +5.93%  '

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '4.257.58'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +17.57%  '

# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '70.892.75'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +6.05%  '

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '3.637.34'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +17.22%  '

# Row 18
$ws.Range("E18").Value = '  +1.93%  '

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '7.55'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +7.35%  '

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '521.41'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +8.44%  '

# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '16.93'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +1.57%  '

# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '9.29'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +18.79%  '

# Row 23
$ws.Range("E23").Value = '  +7.75%  '

# Row 24
$ws.Range("B24").Value = 'Fetch.AI'
$ws.Range("C24").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '2.55'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +13.81%  '

# Row 25
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '88.83'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +6.17%  '

# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '13.44'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +6.68%  '

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '11.09'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +9.93%  '

# Row 28
$ws.Range("E28").Value = '  -0.07%  '

# Row 29
$ws.Range("E29").Value = '  +11.73%  '

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '8.17'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +3.62%  '

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '2.90'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +11.28%  '

# Row 32
$ws.Range("E32").Value = '  +12.96%  '

# Row 33
$ws.Range("E33").Value = '  +17.22%  '

# Row 34
$ws.Range("E34").Value = '  +3.98%  '

# Row 35
$ws.Range("E35").Value = '  +0.01%  '

# Row 36
$ws.Range("E36").Value = '  +9.38%  '

# Row 37
$ws.Range("E37").Value = '  +8.79%  '

# Row 38
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.350'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +12.43%  '

# Row 39
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '2.20'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +9.78%  '

# Row 40
$ws.Range("E40").Value = '  +6.86%  '

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '51.46'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +4.96%  '

# Row 42
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '45.85'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -5.21%  '

# Row 43
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '427.25'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +13.95%  '

# Row 44
$ws.Range("B44").Value = 'Cosmos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '8.83'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +5.91%  '

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '3.113.98'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +11.15%  '

# Row 46
$ws.Range("E46").Value = '  +3.94%  '

# Row 47
$ws.Range("E47").Value = '  +7.70%  '

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '28.41'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +10.89%  '

# Row 49
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '140.57'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +3.85%  '

# Row 51
$ws.Range("E51").Value = '  +10.32%  '
